$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header text updates (Volume/Number and Report Covering the Week)
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/13/2025  Through  1/19/2025"

# ---------------------------------------------------------------------
# 2) Column width updates (columns E and H widen to fit new values)
# ---------------------------------------------------------------------
$ws.Range("E1").EntireColumn.ColumnWidth = 8.0
$ws.Range("H1").EntireColumn.ColumnWidth = 8.0

# ---------------------------------------------------------------------
# 3) Row 16 (Fel. Assault) - "28 Day" block (F,G,H) reset to blank/no-data
# ---------------------------------------------------------------------
$ws.Range("C16").Copy()
$ws.Range("F16").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("C16").Copy()
$ws.Range("F16").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C16").Copy()
$ws.Range("G16").PasteSpecial(-4163)
$ws.Range("C16").Copy()
$ws.Range("G16").PasteSpecial(-4122)

$ws.Range("E16").Copy()
$ws.Range("H16").PasteSpecial(-4163)
$ws.Range("E16").Copy()
$ws.Range("H16").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4) Row 19 (Gr. Larceny) - new data for Week-to-Date / 28-Day / YTD
# ---------------------------------------------------------------------
# D19 <- style/value from G19 (currently 1, style 14), then set to 1
$ws.Range("G19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("G19").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D19").Value = 1

# E19 <- style/value from H19 (currently -100, style 15), then set to -100
$ws.Range("H19").Copy()
$ws.Range("E19").PasteSpecial(-4163)
$ws.Range("H19").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E19").Value = -100

# J19 <- style/value from G19 (still 1, style 14 - before G19 itself changes), then set to 1
$ws.Range("G19").Copy()
$ws.Range("J19").PasteSpecial(-4163)
$ws.Range("G19").Copy()
$ws.Range("J19").PasteSpecial(-4122)
$ws.Range("J19").Value = 1

# K19 <- style/value from H19 (still -100, style 15), then set to -100
$ws.Range("H19").Copy()
$ws.Range("K19").PasteSpecial(-4163)
$ws.Range("H19").Copy()
$ws.Range("K19").PasteSpecial(-4122)
$ws.Range("K19").Value = -100

# G19 keeps style 14, numeric value changes from 1 -> 2
$ws.Range("G19").Value = 2

# H19 and I19 are unchanged by the diff.

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5) Row 21 (TOTAL) - new data for Week-to-Date / 28-Day / YTD
# ---------------------------------------------------------------------
# D21 <- style/value from F21 (currently 1, style 18), then set to 1
$ws.Range("F21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("F21").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("D21").Value = 1

# E21 <- style/value from K21 (currently -100, style 19), then set to -100
$ws.Range("K21").Copy()
$ws.Range("E21").PasteSpecial(-4163)
$ws.Range("K21").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("E21").Value = -100

# F21 <- style/value from C21 (text "0", style 17)
$ws.Range("C21").Copy()
$ws.Range("F21").PasteSpecial(-4163)
$ws.Range("C21").Copy()
$ws.Range("F21").PasteSpecial(-4122)

# G21 keeps style 18, numeric value changes from 5 -> 4
$ws.Range("G21").Value = 4

# H21 keeps style 19, numeric value changes from -80 -> -100
$ws.Range("H21").Value = -100

# J21 keeps style 18, numeric value changes from 1 -> 2
$ws.Range("J21").Value = 2

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 6) Row 27 (UCR Rape*) - new data for Week-to-Date / 28-Day / YTD
# ---------------------------------------------------------------------
# D27 <- style/value from J27 (currently 1, style 14), then set to 1
$ws.Range("J27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("J27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1

# E27 <- style/value from H27 (currently -100, style 15), then set to -100
$ws.Range("H27").Copy()
$ws.Range("E27").PasteSpecial(-4163)
$ws.Range("H27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100

# G27 keeps style 14, numeric value changes from 2 -> 3
$ws.Range("G27").Value = 3

# J27 keeps style 14, numeric value changes from 1 -> 2
$ws.Range("J27").Value = 2

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 7) Row 28 (Other Sex Crimes) - F28 reset to blank/no-data
# ---------------------------------------------------------------------
$ws.Range("C28").Copy()
$ws.Range("F28").PasteSpecial(-4163)
$ws.Range("C28").Copy()
$ws.Range("F28").PasteSpecial(-4122)

$excel.CutCopyMode = 0
